$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Add_Devices_LoopA")
$wsB = $wb.Worksheets.Item("Add_Devices_LoopB")

# Update the test case name string in sheet Add_Devices_LoopA (B2)
$wsA.Range("B2").Value = "verifyTripCurrentCalculationForFIMLoopFC"

# Update numeric values
$wsA.Range("J9").Value = 6
$wsA.Range("J10").Value = 6
$wsB.Range("J10").Value = 6

# Update selections / active sheet to match diff
$wsB.Range("A10").Select()
$wsA.Range("B2").Select()
$wsA.Activate()
